$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Capture handles to the existing conditional-formatting rules
#    BEFORE inserting columns, using single-cell probes that are
#    unambiguous in the ORIGINAL (pre-insert) layout.
# ------------------------------------------------------------------
$ruleNotBetweenG = $ws.Range("G2").FormatConditions.Item(1)   # old G2:G8 notBetween
$ruleNotBetweenJ = $ws.Range("J2").FormatConditions.Item(1)   # old J2:J8 notBetween
$ruleNotBetweenK = $ws.Range("K2").FormatConditions.Item(1)   # old K2:K8 notBetween

$ruleTrue = $null
$ruleFalse = $null
$hRules = $ws.Range("H2").FormatConditions
for ($i = 1; $i -le $hRules.Count; $i++) {
    $it = $hRules.Item($i)
    if ($it.Formula1 -eq "=TRUE") { $ruleTrue = $it }
    if ($it.Formula1 -eq "=FALSE") { $ruleFalse = $it }
}

# ------------------------------------------------------------------
# 2. Insert two new columns (G:H) to make room for the new
#    MPSP_new / GWP_new data; this pushes the old G..L columns to I..N
#    and inherits number formatting/borders from the neighboring column.
# ------------------------------------------------------------------
$ws.Columns("G:H").Insert()

# ------------------------------------------------------------------
# 3. Rename the existing "wwt" columns to "exist" and add headers for
#    the two new columns.
# ------------------------------------------------------------------
$ws.Range("E1").Value = "MPSP_exist"
$ws.Range("F1").Value = "GWP_exist"
$ws.Range("G1").Value = "MPSP_new"
$ws.Range("H1").Value = "GWP_new"

# ------------------------------------------------------------------
# 4. Fill in the new MPSP_new / GWP_new values (columns G:H).
# ------------------------------------------------------------------
$ws.Range("G2").Value2 = 1.44619981691902
$ws.Range("H2").Value2 = 2.6441211717736399

$ws.Range("G3").Value2 = 2.0207467422379302
$ws.Range("H3").Value2 = -1.7700895212399199

$ws.Range("G4").Value2 = 1.6392340672605901
$ws.Range("H4").Value2 = -13.3946003253452

$ws.Range("G5").Value2 = 1.96345339780636
$ws.Range("H5").Value2 = -0.30890062078379998

$ws.Range("G6").Value2 = 2.1099166755697598
$ws.Range("H6").Value2 = 1.1273919300223101

$ws.Range("G7").Value2 = 1.7523427986257001
$ws.Range("H7").Value2 = 0.68677879509477802

$ws.Range("G8").Value2 = 1.3029467468828
$ws.Range("H8").Value2 = 4.2002867619759003

# ------------------------------------------------------------------
# 5. Update the cornstover row (row 7) MPSP_original / GWP_original
#    values after the bug fix.
# ------------------------------------------------------------------
$ws.Range("C7").Value2 = 1.9996562891243299
$ws.Range("D7").Value2 = 1.4157372841731599

# ------------------------------------------------------------------
# 6. Move the conditional-formatting rules onto their new ranges.
# ------------------------------------------------------------------
$ruleNotBetweenG.ModifyAppliesToRange($ws.Range("I2:I8"))
$ruleNotBetweenJ.ModifyAppliesToRange($ws.Range("L2:L8"))
$ruleNotBetweenK.ModifyAppliesToRange($ws.Range("M2:M8"))
$ruleTrue.ModifyAppliesToRange($ws.Range("J2:K8"))
$ruleFalse.ModifyAppliesToRange($ws.Range("J2:K8"))

# extend the TRUE/FALSE rules to also cover the new consistency column N
$nTrue = $ws.Range("N2:N8").FormatConditions.Add(1, 3, "=TRUE")
$nFalse = $ws.Range("N2:N8").FormatConditions.Add(1, 3, "=FALSE")
$nTrue.Interior.ColorIndex = $ruleTrue.Interior.ColorIndex
$nTrue.Font.ColorIndex = $ruleTrue.Font.ColorIndex
$nFalse.Interior.ColorIndex = $ruleFalse.Interior.ColorIndex
$nFalse.Font.ColorIndex = $ruleFalse.Font.ColorIndex

# ------------------------------------------------------------------
# 7. Move the comment from (old) F7 to A7.
# ------------------------------------------------------------------
$commentText = $ws.Range("F7").Comment.Text()
$ws.Range("F7").Comment.Delete()
$ws.Range("A7").AddComment($commentText)

# ------------------------------------------------------------------
# 8. Update the workbook view window position & selection.
# ------------------------------------------------------------------
$excel.ActiveWindow.Left = 4880
$excel.ActiveWindow.Top = 10420
$ws.Range("H8").Select()
